# Rename the existing sheet from "Sheet1" to "Data"
$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item(1)
$dataSheet.Name = "Data"

# Add a new worksheet "FoF" right after the "Data" sheet
$fofSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $dataSheet)
$fofSheet.Name = "FoF"

# Header
$fofSheet.Range("B1").Value = "FoF"

# Stats block
$fofSheet.Range("A2").Value = "Return"
$fofSheet.Range("B2").Value = 0.049341

$fofSheet.Range("A3").Value = "M squared"
$fofSheet.Range("B3").Value = 0.049341

$fofSheet.Range("A4").Value = "Volatility"
$fofSheet.Range("B4").Value = 0.056391

$fofSheet.Range("A5").Value = "MDD"
$fofSheet.Range("B5").Value = 0.222035

$fofSheet.Range("A6").Value = "CVaR"
$fofSheet.Range("B6").Value = 0.072119

$fofSheet.Range("A7").Value = "CDaR"
$fofSheet.Range("B7").Value = 0.242259

$fofSheet.Range("A8").Value = "Sharpe"
$fofSheet.Range("B8").Value = 0.122638

$fofSheet.Range("A9").Value = "Calmar"
$fofSheet.Range("B9").Value = 0.031147

$fofSheet.Range("A10").Value = "R squared"
$fofSheet.Range("B10").Value = 1

$fofSheet.Range("A11").Value = "Corr. Stocks"
$fofSheet.Range("B11").Value = 0.647835

$fofSheet.Range("A12").Value = "Corr. Bonds"
$fofSheet.Range("B12").Value = 0.360859

$fofSheet.Range("A13").Value = "Corr. FoF"
$fofSheet.Range("B13").Value = 1

$fofSheet.Range("A14").Value = "Turnover"
$fofSheet.Range("B14").Value = 0

# Make FoF the active sheet/tab, with the same lingering selection as the
# author's session
$fofSheet.Activate()
$fofSheet.Range("E23").Select()
